# Update cryptocurrency price/volume data to latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'23.017.41"
$ws.Cells.Item(2, 5).Value = "'  -2.96%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.598.92"
$ws.Cells.Item(3, 5).Value = "'  -1.83%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "'  -0.21%  "

# Row 5
$ws.Cells.Item(5, 5).Value = "'  -0.15%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'301.74"
$ws.Cells.Item(6, 5).Value = "'  -1.89%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.3783"
$ws.Cells.Item(7, 5).Value = "'  -1.16%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.3642"
$ws.Cells.Item(8, 5).Value = "'  -3.93%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'49.89"
$ws.Cells.Item(9, 5).Value = "'  -1.25%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'1.255"
$ws.Cells.Item(10, 5).Value = "'  -4.32%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "'BinanceUSD"
$ws.Cells.Item(11, 3).Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(11, 4).Value = "'1.002"
$ws.Cells.Item(11, 5).Value = "'  -0.25%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "'Dogecoin"
$ws.Cells.Item(12, 3).Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(12, 4).Value = "'0.08143"
$ws.Cells.Item(12, 5).Value = "'  -2.06%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'22.47"
$ws.Cells.Item(13, 5).Value = "'  -4.27%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'6.591"
$ws.Cells.Item(14, 5).Value = "'  -4.39%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.00001246"
$ws.Cells.Item(15, 5).Value = "'  -3.01%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'7.343"
$ws.Cells.Item(16, 5).Value = "'  -4.51%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'1.601.90"
$ws.Cells.Item(17, 5).Value = "'  -1.83%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'91.91"
$ws.Cells.Item(18, 5).Value = "'  -1.17%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'0.06821"
$ws.Cells.Item(19, 5).Value = "'  -1.46%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'18.20"
$ws.Cells.Item(20, 5).Value = "'  -5.43%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'6.529"
$ws.Cells.Item(21, 5).Value = "'  -4.05%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "'Dai"
$ws.Cells.Item(22, 3).Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(22, 4).Value = "'1.002"
$ws.Cells.Item(22, 5).Value = "'  -0.02%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "'Cosmos"
$ws.Cells.Item(23, 3).Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(23, 4).Value = "'13.05"
$ws.Cells.Item(23, 5).Value = "'  -2.73%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "'WrappedBTC"
$ws.Cells.Item(24, 3).Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(24, 4).Value = "'23.011.44"
$ws.Cells.Item(24, 5).Value = "'  -2.98%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "'Toncoin"
$ws.Cells.Item(25, 3).Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(25, 4).Value = "'2.364"
$ws.Cells.Item(25, 5).Value = "'  -2.10%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "'LidoDAOToken"
$ws.Cells.Item(26, 3).Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(26, 4).Value = "'2.809"
$ws.Cells.Item(26, 5).Value = "'  -0.71%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "'EthereumClassic"
$ws.Cells.Item(27, 3).Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).Value = "'21.05"
$ws.Cells.Item(27, 5).Value = "'  -3.11%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "'Monero"
$ws.Cells.Item(28, 3).Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(28, 4).Value = "'150.22"
$ws.Cells.Item(28, 5).Value = "'  -1.16%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "'HuobiToken"
$ws.Cells.Item(29, 3).Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(29, 4).Value = "'5.232"
$ws.Cells.Item(29, 5).Value = "'  -3.59%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "'BitcoinCash"
$ws.Cells.Item(30, 3).Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(30, 4).Value = "'134.45"
$ws.Cells.Item(30, 5).Value = "'  -1.25%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "'WEMIXTOKEN"
$ws.Cells.Item(31, 3).Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(31, 4).Value = "'2.345"
$ws.Cells.Item(31, 5).Value = "'  -5.34%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "'Filecoin"
$ws.Cells.Item(32, 3).Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).Value = "'6.821"
$ws.Cells.Item(32, 5).Value = "'  -13.86%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "'WrappedliquidstakedEther2.0"
$ws.Cells.Item(33, 3).Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(33, 4).Value = "'1.779.03"
$ws.Cells.Item(33, 5).Value = "'  -1.83%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "'ImmutableX"
$ws.Cells.Item(34, 3).Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).Value = "'0.9616"
$ws.Cells.Item(34, 5).Value = "'  -1.51%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "'Hedera"
$ws.Cells.Item(35, 3).Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(35, 4).Value = "'0.07569"
$ws.Cells.Item(35, 5).Value = "'  -3.24%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "'FraxShare"
$ws.Cells.Item(36, 3).Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(36, 4).Value = "'10.29"
$ws.Cells.Item(36, 5).Value = "'  +0.34%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "'InternetComputer(DFINITY)"
$ws.Cells.Item(37, 3).Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(37, 4).Value = "'6.241"
$ws.Cells.Item(37, 5).Value = "'  -3.81%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "'VeChain"
$ws.Cells.Item(38, 3).Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.02706"
$ws.Cells.Item(38, 5).Value = "'  -5.33%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "'Algorand"
$ws.Cells.Item(39, 3).Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(39, 4).Value = "'0.2524"
$ws.Cells.Item(39, 5).Value = "'  -3.97%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "'Stellar"
$ws.Cells.Item(40, 3).Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(40, 4).Value = "'0.08848"
$ws.Cells.Item(40, 5).Value = "'  -1.71%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "'TrustWalletToken"
$ws.Cells.Item(41, 3).Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(41, 4).Value = "'1.363"
$ws.Cells.Item(41, 5).Value = "'  -3.15%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "'TheSandbox"
$ws.Cells.Item(42, 3).Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(42, 4).Value = "'0.7025"
$ws.Cells.Item(42, 5).Value = "'  -5.11%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "'Aptos"
$ws.Cells.Item(43, 3).Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(43, 4).Value = "'12.38"
$ws.Cells.Item(43, 5).Value = "'  -5.64%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "'EnergySwap"
$ws.Cells.Item(44, 3).Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).Value = "'15.24"
$ws.Cells.Item(44, 5).Value = "'  -6.19%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "'Decentraland"
$ws.Cells.Item(45, 3).Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(45, 4).Value = "'0.6614"
$ws.Cells.Item(45, 5).Value = "'  -2.84%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "'Frax"
$ws.Cells.Item(46, 3).Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(46, 4).Value = "'0.9996"
$ws.Cells.Item(46, 5).Value = "'  -0.15%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "'PancakeSwap"
$ws.Cells.Item(47, 3).Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(47, 4).Value = "'3.994"
$ws.Cells.Item(47, 5).Value = "'  -1.18%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'2.283"
$ws.Cells.Item(48, 5).Value = "'  -4.49%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "'Quant"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).Value = "'132.48"
$ws.Cells.Item(49, 5).Value = "'  -0.59%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "'Cronos"
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.07900"
$ws.Cells.Item(50, 5).Value = "'  -3.27%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "'Flow"
$ws.Cells.Item(51, 3).Value = "'https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Cells.Item(51, 4).Value = "'1.219"
$ws.Cells.Item(51, 5).Value = "'  +1.17%  "
